$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, and whether it must be
# forced to Text so Excel does not silently reinterpret it as a number
# (e.g. "1.00" -> 1, or "213.79" -> numeric 213.79 instead of the literal string).
$updates = @(
    @{ Cell = "D2"; Value = '27.423.45'; ForceText = $true }
    @{ Cell = "E2"; Value = '  -2.09%  '; ForceText = $false }
    @{ Cell = "D3"; Value = '1.653.63'; ForceText = $true }
    @{ Cell = "E3"; Value = '  -1.96%  '; ForceText = $false }
    @{ Cell = "D4"; Value = '0.999'; ForceText = $true }
    @{ Cell = "E4"; Value = '  -0.17%  '; ForceText = $false }
    @{ Cell = "D5"; Value = '213.79'; ForceText = $true }
    @{ Cell = "E5"; Value = '  -1.37%  '; ForceText = $false }
    @{ Cell = "D6"; Value = '0.509'; ForceText = $true }
    @{ Cell = "E6"; Value = '  -1.86%  '; ForceText = $false }
    @{ Cell = "D7"; Value = '0.999'; ForceText = $true }
    @{ Cell = "E7"; Value = '  -0.24%  '; ForceText = $false }
    @{ Cell = "D8"; Value = '24.26'; ForceText = $true }
    @{ Cell = "E8"; Value = '  +1.06%  '; ForceText = $false }
    @{ Cell = "E9"; Value = '  -0.98%  '; ForceText = $false }
    @{ Cell = "E10"; Value = '  -1.57%  '; ForceText = $false }
    @{ Cell = "D11"; Value = '0.0878'; ForceText = $true }
    @{ Cell = "E11"; Value = '  -0.64%  '; ForceText = $false }
    @{ Cell = "D12"; Value = '1.884.78'; ForceText = $true }
    @{ Cell = "E12"; Value = '  -2.12%  '; ForceText = $false }
    @{ Cell = "D13"; Value = '1.649.38'; ForceText = $true }
    @{ Cell = "E13"; Value = '  -2.24%  '; ForceText = $false }
    @{ Cell = "D14"; Value = '0.575'; ForceText = $true }
    @{ Cell = "E14"; Value = '  +2.89%  '; ForceText = $false }
    @{ Cell = "E15"; Value = '  -2.08%  '; ForceText = $false }
    @{ Cell = "D16"; Value = '65.94'; ForceText = $true }
    @{ Cell = "E16"; Value = '  -1.33%  '; ForceText = $false }
    @{ Cell = "D17"; Value = '27.415.73'; ForceText = $true }
    @{ Cell = "E17"; Value = '  -1.96%  '; ForceText = $false }
    @{ Cell = "D18"; Value = '234.37'; ForceText = $true }
    @{ Cell = "E18"; Value = '  -6.26%  '; ForceText = $false }
    @{ Cell = "D19"; Value = '0.0₃0728'; ForceText = $true }
    @{ Cell = "E19"; Value = '  -1.93%  '; ForceText = $false }
    @{ Cell = "D20"; Value = '7.46'; ForceText = $true }
    @{ Cell = "E20"; Value = '  -2.88%  '; ForceText = $false }
    @{ Cell = "D21"; Value = '1.00'; ForceText = $true }
    @{ Cell = "E21"; Value = '  -0.10%  '; ForceText = $false }
    @{ Cell = "E22"; Value = '  -2.55%  '; ForceText = $false }
    @{ Cell = "D23"; Value = '9.32'; ForceText = $true }
    @{ Cell = "E23"; Value = '  -2.28%  '; ForceText = $false }
    @{ Cell = "E24"; Value = '  -1.08%  '; ForceText = $false }
    @{ Cell = "D25"; Value = '147.06'; ForceText = $true }
    @{ Cell = "E25"; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = "D26"; Value = '7.22'; ForceText = $true }
    @{ Cell = "E26"; Value = '  -1.54%  '; ForceText = $false }
    @{ Cell = "D27"; Value = '16.02'; ForceText = $true }
    @{ Cell = "E27"; Value = '  -2.71%  '; ForceText = $false }
    @{ Cell = "D28"; Value = '0.999'; ForceText = $true }
    @{ Cell = "E28"; Value = '  -0.19%  '; ForceText = $false }
    @{ Cell = "E29"; Value = '  -2.09%  '; ForceText = $false }
    @{ Cell = "E30"; Value = '  -1.14%  '; ForceText = $false }
    @{ Cell = "D31"; Value = '1.20'; ForceText = $true }
    @{ Cell = "E31"; Value = '  -5.64%  '; ForceText = $false }
    @{ Cell = "E32"; Value = '  -2.22%  '; ForceText = $false }
    @{ Cell = "D33"; Value = '1.463.32'; ForceText = $true }
    @{ Cell = "E33"; Value = '  +2.45%  '; ForceText = $false }
    @{ Cell = "D34"; Value = '3.11'; ForceText = $true }
    @{ Cell = "E34"; Value = '  -2.21%  '; ForceText = $false }
    @{ Cell = "E35"; Value = '  -3.96%  '; ForceText = $false }
    @{ Cell = "E36"; Value = '  -0.77%  '; ForceText = $false }
    @{ Cell = "D37"; Value = '0.912'; ForceText = $true }
    @{ Cell = "E37"; Value = '  -3.26%  '; ForceText = $false }
    @{ Cell = "E38"; Value = '  -3.39%  '; ForceText = $false }
    @{ Cell = "E39"; Value = '  -1.32%  '; ForceText = $false }
    @{ Cell = "E40"; Value = '  +0.06%  '; ForceText = $false }
    @{ Cell = "D41"; Value = '0.999'; ForceText = $true }
    @{ Cell = "E41"; Value = '  -0.18%  '; ForceText = $false }
    @{ Cell = "D42"; Value = '5.47'; ForceText = $true }
    @{ Cell = "E42"; Value = '  -0.62%  '; ForceText = $false }
    @{ Cell = "E43"; Value = '  -5.50%  '; ForceText = $false }
    @{ Cell = "E44"; Value = '  -0.75%  '; ForceText = $false }
    @{ Cell = "E45"; Value = '  -1.42%  '; ForceText = $false }
    @{ Cell = "D46"; Value = '1.793.76'; ForceText = $true }
    @{ Cell = "E46"; Value = '  -2.17%  '; ForceText = $false }
    @{ Cell = "E47"; Value = '  +0.41%  '; ForceText = $false }
    @{ Cell = "D48"; Value = '88.38'; ForceText = $true }
    @{ Cell = "E48"; Value = '  -1.00%  '; ForceText = $false }
    @{ Cell = "E49"; Value = '  -3.68%  '; ForceText = $false }
    @{ Cell = "E50"; Value = '  -1.55%  '; ForceText = $false }
    @{ Cell = "D51"; Value = '7.78'; ForceText = $true }
    @{ Cell = "E51"; Value = '  -2.01%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Apply a text number-format so the literal digits/dots are kept as-is,
        # then clear the format again so the cell keeps its original (default) style.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
